$wb = $excel.ActiveWorkbook

# --- Rename existing sheets ---
$wb.Worksheets.Item("Sheet1").Name = "original_idea"
$wb.Worksheets.Item("Sheet2").Name = "table"

$table = $wb.Worksheets.Item("table")

# --- Update the "table" sheet view: drop tab selection, change selection ---
$table.Activate()
$table.Range("D3:H4").Select()

# --- Add the new "mapping" sheet after "table" ---
$mapping = $wb.Worksheets.Add($null, $table)
$mapping.Name = "mapping"

# --- Fill in the mapping data (write order controls shared-string order) ---
$mapping.Cells.Item(1,2).Value = "top_left"
$mapping.Cells.Item(1,3).Value = "bottom_right"

$mapping.Cells.Item(2,1).Value = "TR - P1 - R1 - C1 "
$mapping.Cells.Item(2,2).Value = " TL - P3 - R2 - C3 "
$mapping.Cells.Item(2,3).Value = " BR - P2 - R3 - C2"

$mapping.Cells.Item(3,1).Value = "TR - P3 - R2 - C1 "
$mapping.Cells.Item(3,2).Value = " TL - P1 - R1 - C3 "
$mapping.Cells.Item(3,3).Value = " BR - P4 - R4 - C2"

$mapping.Cells.Item(4,1).Value = "TR - P4 - R2 - C2 "
$mapping.Cells.Item(4,2).Value = " TL - P2 - R1 - C4 "
$mapping.Cells.Item(4,3).Value = " BR - P3 - R4 - C1"

$mapping.Cells.Item(5,1).Value = "TR - P2 - R1 - C2 "
$mapping.Cells.Item(5,2).Value = " TL - P4 - R2 - C4 "
$mapping.Cells.Item(5,3).Value = " BR - P1 - R3 - C1"

$mapping.Cells.Item(1,1).Value = "top_right"

# --- Column widths (best-fit to content) ---
$mapping.Columns.Item(1).AutoFit()
$mapping.Columns.Item(2).AutoFit()
$mapping.Columns.Item(3).AutoFit()

# --- View: zoom + selection on the mapping sheet ---
$mapping.Activate()
$excel.ActiveWindow.Zoom = 172
$mapping.Range("D2").Select()

Write-Host "done"
